# Commit: "fix bug exeded requeste in google drive"
# The sheet "Hoja1" header date (A1) is bumped by one day, and the price
# cell for the "TER-10" row (D44) is corrected to a new numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: date serial 45310 -> 45311 (2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# D44: price value 43783.243 -> 17516.625
$ws.Range("D44").Value = 17516.625
